$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.48"
$ws.Range("E2").Value = "'3.98%"
$ws.Range("D3").Value = "'35.64"
$ws.Range("E3").Value = "'14.08%"
$ws.Range("D4").Value = "'5.094"
$ws.Range("E4").Value = "'2.65%"
$ws.Range("D5").Value = "'0.07830"
$ws.Range("E5").Value = "'4.45%"
$ws.Range("D6").Value = "'2.243"
$ws.Range("E6").Value = "'-1.32%"
$ws.Range("D7").Value = "'8.124"
$ws.Range("E7").Value = "'4.15%"
$ws.Range("D8").Value = "'4.008"
$ws.Range("E8").Value = "'6.35%"
$ws.Range("D9").Value = "'0.9269"
$ws.Range("E9").Value = "'0.68%"
$ws.Range("D10").Value = "'0.09695"
$ws.Range("E10").Value = "'2.99%"
$ws.Range("D11").Value = "'0.1820"
$ws.Range("E11").Value = "'4.76%"
$ws.Range("D12").Value = "'0.08709"
$ws.Range("E12").Value = "'4.01%"
$ws.Range("E13").Value = "'4.43%"
$ws.Range("D14").Value = "'0.09934"
$ws.Range("E14").Value = "'-0.09%"
$ws.Range("D15").Value = "'0.001480"
$ws.Range("E15").Value = "'-1.37%"
$ws.Range("D16").Value = "'0.005706"
$ws.Range("E16").Value = "'-1.57%"
$ws.Range("D17").Value = "'3.485"
$ws.Range("E17").Value = "'0.28%"
$ws.Range("E18").Value = "'-2.18%"
$ws.Range("D19").Value = "'0.3459"
$ws.Range("E19").Value = "'3.45%"
$ws.Range("E20").Value = "'0.48%"
$ws.Range("D21").Value = "'4.544"
$ws.Range("E21").Value = "'10.32%"
$ws.Range("D23").Value = "'0.04680"
$ws.Range("E23").Value = "'3.21%"
$ws.Range("D24").Value = "'0.001245"
$ws.Range("E24").Value = "'2.16%"
$ws.Range("D25").Value = "'0.004540"
$ws.Range("E25").Value = "'5.30%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'0.10%"
$ws.Range("D27").Value = "'0.0002699"
$ws.Range("E27").Value = "'-20.40%"
$ws.Range("D39").Value = "'0.01758"
$ws.Range("E39").Value = "'6.96%"
$ws.Range("D40").Value = "'0.04711"
$ws.Range("E40").Value = "'2.50%"
$ws.Range("D41").Value = "'0.007916"
$ws.Range("E41").Value = "'5.89%"
$ws.Range("E42").Value = "'4.60%"
$ws.Range("D43").Value = "'0.008041"
$ws.Range("E43").Value = "'-18.21%"
$ws.Range("D44").Value = "'0.002300"
$ws.Range("E44").Value = "'3.71%"
$ws.Range("D45").Value = "'0.009131"
$ws.Range("E45").Value = "'-3.00%"
$ws.Range("D46").Value = "'0.00006211"
$ws.Range("E46").Value = "'1.83%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.12%"
$ws.Range("D48").Value = "'4.005"
$ws.Range("E48").Value = "'57.01%"
$ws.Range("D49").Value = "'0.002690"
$ws.Range("E49").Value = "'34.65%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'0.12%"
